$wb = $excel.ActiveWorkbook

# --- Sheet "V3": fill in the B-column frame numbers for rows 72-89 ---
# (column D recalculates automatically via the existing shared formula
#  D = IF(C<>"", IF(B<>"", C-B, "-"), "-"))
$wsV3 = $wb.Worksheets.Item("V3")

$wsV3.Range("B72").Value = 46958
$wsV3.Range("B73").Value = 47930
$wsV3.Range("B74").Value = 48273
$wsV3.Range("B75").Value = 49356
$wsV3.Range("B76").Value = 49699
$wsV3.Range("B77").Value = 50839
$wsV3.Range("B78").Value = 51182
$wsV3.Range("B79").Value = 51882
$wsV3.Range("B80").Value = 52225
$wsV3.Range("B81").Value = 53516
$wsV3.Range("B82").Value = 53860
$wsV3.Range("B83").Value = 55029
$wsV3.Range("B84").Value = 55373
$wsV3.Range("B85").Value = 56774
$wsV3.Range("B86").Value = 57118
$wsV3.Range("B87").Value = 57747
$wsV3.Range("B88").Value = 58090
$wsV3.Range("B89").Value = 59219

# Move the view/selection on "V3" down to where the new data was entered
$wsV3.Activate()
$excel.ActiveWindow.ScrollRow = 74
$excel.ActiveWindow.ScrollColumn = 1
$wsV3.Range("B90").Select()

# --- Sheet "V2": just the view/selection moved ---
$wsV2 = $wb.Worksheets.Item("V2")
$wsV2.Activate()
$excel.ActiveWindow.ScrollRow = 61
$excel.ActiveWindow.ScrollColumn = 1
$wsV2.Range("C77").Select()

# Leave "V3" as the active/selected sheet, matching tabSelected="1" on V3
$wsV3.Activate()
